$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "jumlah" (quantity) value in I2 from 2 to 3
$ws.Range("I2").Value = 3

# Move active cell selection to I3 (matches updated selection in sheet view)
$ws.Range("I3").Select()
